$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name in B3
$ws.Range("B3").Value = "Otello Corporation ASA (OB:OTEC)"

# Update numeric values for row 2 and row 3 (same new values in both rows)
$ws.Range("D2").Value = -0.126
$ws.Range("D3").Value = -0.126
$ws.Range("G2").Value = 0.1041580897488678
$ws.Range("G3").Value = 0.1041580897488678
$ws.Range("H2").Value = -0.0168793742280774
$ws.Range("H3").Value = -0.0168793742280774
$ws.Range("I2").Value = -0.03417044051049815
$ws.Range("I3").Value = -0.03417044051049815
$ws.Range("J2").Value = -0.03417044051049815
$ws.Range("J3").Value = -0.03417044051049815
$ws.Range("K2").Value = -19.8
$ws.Range("K3").Value = -19.8
$ws.Range("L2").Value = -0.08151502675998354
$ws.Range("L3").Value = -0.08151502675998354
$ws.Range("M2").Value = 0.6
$ws.Range("M3").Value = 0.6
$ws.Range("N2").Value = 0.001379310344827586
$ws.Range("N3").Value = 0.001379310344827586
$ws.Range("O2").Value = -0.0303030303030303
$ws.Range("O3").Value = -0.0303030303030303
$ws.Range("S2").Value = 0.6
$ws.Range("S3").Value = 0.6
$ws.Range("U2").Value = 33.1
$ws.Range("U3").Value = 33.1
$ws.Range("V2").Value = 0.0760919540229885
$ws.Range("V3").Value = 0.0760919540229885
$ws.Range("W2").Value = -0.05774278215223098
$ws.Range("W3").Value = -0.05774278215223098
$ws.Range("X2").Value = 0.05882414688956489
$ws.Range("X3").Value = 0.05882414688956489
$ws.Range("Y2").Value = -0.1165669290417959
$ws.Range("Y3").Value = -0.1165669290417959
$ws.Range("Z2").Value = 2.222323879231473
$ws.Range("Z3").Value = 2.222323879231473
$ws.Range("AA2").Value = -0.07593778591033853
$ws.Range("AA3").Value = -0.07593778591033853
$ws.Range("AB2").Value = 0.05598800489339737
$ws.Range("AB3").Value = 0.05598800489339737
$ws.Range("AC2").Value = -0.1319257908037359
$ws.Range("AC3").Value = -0.1319257908037359
$ws.Range("AD2").Value = 34.4
$ws.Range("AD3").Value = 34.4
$ws.Range("AE2").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF2").Value = 34.4
$ws.Range("AF3").Value = 34.4
$ws.Range("AG2").Value = 1.299999999999997
$ws.Range("AG3").Value = 1.299999999999997
$ws.Range("AH2").Value = 0.07328504473796336
$ws.Range("AH3").Value = 0.07328504473796336
$ws.Range("AI2").Value = 0.1029940119760479
$ws.Range("AI3").Value = 0.1029940119760479
$ws.Range("AJ2").Value = 0.00297960119184047
$ws.Range("AJ3").Value = 0.00297960119184047
$ws.Range("AK2").Value = 0.004320372216683273
$ws.Range("AK3").Value = 0.004320372216683273
$ws.Range("AL2").Value = 1.7
$ws.Range("AL3").Value = 1.7
$ws.Range("AM2").Value = 0.7
$ws.Range("AM3").Value = 0.7
$ws.Range("AN2").Value = -24.57142857142857
$ws.Range("AN3").Value = -24.57142857142857
$ws.Range("AO2").Value = -4.882352941176471
$ws.Range("AO3").Value = -4.882352941176471
$ws.Range("AP2").Value = -0.9285714285714266
$ws.Range("AP3").Value = -0.9285714285714266
$ws.Range("AQ2").Value = -11.85714285714286
$ws.Range("AQ3").Value = -11.85714285714286
